$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 454.87097
$ws.Range("I15").Value = 454.87097
$ws.Range("K15").Value = 1364.61291
$ws.Range("M15").Value = -1195.61291

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1501.7179
$ws.Range("J17").Value = 1501.7179
$ws.Range("L17").Value = 4505.153700000001
$ws.Range("N17").Value = -4841.153700000001

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 33039.17
$ws.Range("I62").Value = 104737
$ws.Range("K62").Value = 104737
$ws.Range("M62").Value = -104113

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 33039.17
$ws.Range("I65").Value = 104737
$ws.Range("K65").Value = 523685
$ws.Range("M65").Value = -520565

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3099.95
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 5000
$ws.Range("N76").Value = -5630

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3099.95
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 5000
$ws.Range("N79").Value = -7184

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4996
$ws.Range("J116").Value = 4666.6665
$ws.Range("L116").Value = 4666.6665
$ws.Range("N116").Value = -11550.6665

# ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 399999
$ws.Range("J136").Value = 399999
$ws.Range("L136").Value = 399999
$ws.Range("N136").Value = -410199

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2597.8474
$ws.Range("I138").Value = 873.7368
$ws.Range("J138").Value = 3416.8
$ws.Range("K138").Value = 2621.2104
$ws.Range("L138").Value = 10250.4
$ws.Range("M138").Value = 2518.7896
$ws.Range("N138").Value = -20530.4

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2648
$ws.Range("I141").Value = 2537
$ws.Range("K141").Value = 7611
$ws.Range("M141").Value = -2431

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 323.33334
$ws.Range("I5").Value = 110
$ws.Range("K5").Value = 110
$ws.Range("M5").Value = 2

# ARM row 31
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N31").ClearContents()
$ws.Range("H31").Value = 4490.3335
$ws.Range("I31").Value = 4490.3335
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 4490.3335
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -4196.3335

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2473.4736
$ws.Range("I32").Value = 2357.574
$ws.Range("J32").Value = 4559.6665
$ws.Range("K32").Value = 2357.574
$ws.Range("L32").Value = 4559.6665
$ws.Range("M32").Value = -2070.574
$ws.Range("N32").Value = -5133.6665

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2642.5715
$ws.Range("I74").Value = 2816.3333
$ws.Range("J74").Value = 1600
$ws.Range("K74").Value = 2816.3333
$ws.Range("L74").Value = 1600
$ws.Range("M74").Value = -1942.3333
$ws.Range("N74").Value = -3348

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2642.5715
$ws.Range("I77").Value = 2816.3333
$ws.Range("J77").Value = 1600
$ws.Range("K77").Value = 14081.6665
$ws.Range("L77").Value = 8000
$ws.Range("M77").Value = -9713.666499999999
$ws.Range("N77").Value = -16736

# ARM row 92
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1943.2916
$ws.Range("I122").Value = 1775.7368
$ws.Range("J122").Value = 2580
$ws.Range("K122").Value = 5327.2104
$ws.Range("L122").Value = 7740
$ws.Range("M122").Value = -2877.2104
$ws.Range("N122").Value = -12640

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 323.33334
$ws.Range("I4").Value = 110
$ws.Range("K4").Value = 110
$ws.Range("M4").Value = 5

# BSM row 7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 728
$ws.Range("I7").Value = 575
$ws.Range("J7").Value = 830
$ws.Range("K7").Value = 575
$ws.Range("L7").Value = 830
$ws.Range("M7").Value = -462
$ws.Range("N7").Value = -1056

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2029.3077
$ws.Range("I86").Value = 2129.3333
$ws.Range("J86").Value = 829
$ws.Range("K86").Value = 2129.3333
$ws.Range("L86").Value = 829
$ws.Range("M86").Value = -1006.3333
$ws.Range("N86").Value = -3075

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2029.3077
$ws.Range("I89").Value = 2129.3333
$ws.Range("J89").Value = 829
$ws.Range("K89").Value = 10646.6665
$ws.Range("L89").Value = 4145
$ws.Range("M89").Value = -5030.666499999999
$ws.Range("N89").Value = -15377

# CRP row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2031.6666
$ws.Range("I12").Value = 1100
$ws.Range("K12").Value = 1100
$ws.Range("M12").Value = -930

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1716.8235
$ws.Range("I22").Value = 1613.5714
$ws.Range("J22").Value = 1789.1
$ws.Range("K22").Value = 1613.5714
$ws.Range("L22").Value = 1789.1
$ws.Range("M22").Value = -1263.5714
$ws.Range("N22").Value = -2489.1

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 28149.117
$ws.Range("J41").Value = 38164.668
$ws.Range("L41").Value = 38164.668
$ws.Range("N41").Value = -39020.668

# CRP row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 15487
$ws.Range("J88").Value = 15487
$ws.Range("L88").Value = 15487
$ws.Range("N88").Value = -16299

# CRP row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 15487
$ws.Range("J91").Value = 15487
$ws.Range("L91").Value = 15487
$ws.Range("N91").Value = -18295

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4815
$ws.Range("I99").Value = 3518.75
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 3518.75
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -2020.75
$ws.Range("N99").Value = -12996

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4815
$ws.Range("I126").Value = 3518.75
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 10556.25
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -8086.25
$ws.Range("N126").Value = -34940

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2709.0293
$ws.Range("I132").Value = 2803.516
$ws.Range("J132").Value = 1732.6666
$ws.Range("K132").Value = 8410.548000000001
$ws.Range("L132").Value = 5197.9998
$ws.Range("M132").Value = -5880.548000000001
$ws.Range("N132").Value = -10257.9998

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 333.6154
$ws.Range("J12").Value = 404.77777
$ws.Range("L12").Value = 1214.33331
$ws.Range("N12").Value = -1560.33331

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2037.5
$ws.Range("I140").Value = 1383.3334
$ws.Range("K140").Value = 4150.0002
$ws.Range("M140").Value = 1029.9998

# GSM row 22
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 8779.777
$ws.Range("J22").Value = 12004.5
$ws.Range("L22").Value = 12004.5
$ws.Range("N22").Value = -13062.5

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6937.4
$ws.Range("I70").Value = 6647.3335
$ws.Range("K70").Value = 6647.3335
$ws.Range("M70").Value = -6377.3335

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6937.4
$ws.Range("I73").Value = 6647.3335
$ws.Range("K73").Value = 6647.3335
$ws.Range("M73").Value = -5711.3335

# GSM row 96
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 25000
$ws.Range("J96").Value = 25000
$ws.Range("L96").Value = 25000
$ws.Range("N96").Value = -30492

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1180.6
$ws.Range("I97").Value = 1066
$ws.Range("K97").Value = 1066
$ws.Range("M97").Value = -570

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2737.2222
$ws.Range("I126").Value = 2737.2222
$ws.Range("K126").Value = 8211.6666
$ws.Range("M126").Value = -5741.6666

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3589.6843
$ws.Range("I132").Value = 3387.8125
$ws.Range("K132").Value = 10163.4375
$ws.Range("M132").Value = -7633.4375

# GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 59999.25
$ws.Range("J138").Value = 59999.25
$ws.Range("L138").Value = 59999.25
$ws.Range("N138").Value = -70279.25

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4080
$ws.Range("I16").Value = 1250
$ws.Range("K16").Value = 1250
$ws.Range("M16").Value = -1080

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4927.222
$ws.Range("I132").Value = 4293.125
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 12879.375
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -10349.375
$ws.Range("N132").Value = -35060

# WVR row 3
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 551.5
$ws.Range("J3").Value = 1000
$ws.Range("L3").Value = 1000
$ws.Range("N3").Value = -1228

# WVR row 11
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M11").ClearContents()
$ws.Range("H11").Value = 6666.6665
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 6666.6665
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 6666.6665
$ws.Range("N11").Value = -6950.6665

# WVR row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 44979
$ws.Range("J70").Value = 44979
$ws.Range("L70").Value = 44979
$ws.Range("N70").Value = -45609

# WVR row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 44979
$ws.Range("J73").Value = 44979
$ws.Range("L73").Value = 44979
$ws.Range("N73").Value = -47163

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1202.1794
$ws.Range("I132").Value = 1053.6285
$ws.Range("K132").Value = 3160.8855
$ws.Range("M132").Value = -630.8855000000003
